$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 5131.4
$ws.Range("I6").Value = 5131.4
$ws.Range("K6").Value = 15394.2
$ws.Range("M6").Value = -15282.2
$ws.Range("H33").Value = 261.22223
$ws.Range("I33").Value = 145.66667
$ws.Range("J33").Value = 492.33334
$ws.Range("K33").Value = 145.66667
$ws.Range("L33").Value = 492.33334
$ws.Range("M33").Value = 83.33332999999999
$ws.Range("N33").Value = -950.33334
$ws.Range("H55").Value = 213.91667
$ws.Range("I55").Value = 172.14285
$ws.Range("J55").Value = 272.4
$ws.Range("K55").Value = 172.14285
$ws.Range("L55").Value = 272.4
$ws.Range("M55").Value = 41.85714999999999
$ws.Range("N55").Value = -700.4
$ws.Range("H137").Value = 1217.2354
$ws.Range("I137").Value = 1233.0769
$ws.Range("J137").Value = 1165.75
$ws.Range("K137").Value = 3699.2307
$ws.Range("L137").Value = 3497.25
$ws.Range("M137").Value = -1149.2307
$ws.Range("N137").Value = -8597.25
$ws.Range("H138").Value = 1279.13
$ws.Range("I138").Value = 664.8158
$ws.Range("J138").Value = 1655.6451
$ws.Range("K138").Value = 1994.4474
$ws.Range("L138").Value = 4966.9353
$ws.Range("M138").Value = 3145.5526
$ws.Range("N138").Value = -15246.9353

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 445
$ws.Range("I4").Value = 445
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 445
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -329
$ws.Range("N4").ClearContents()
$ws.Range("H9").Value = 10009
$ws.Range("J9").Value = 10009
$ws.Range("L9").Value = 10009
$ws.Range("N9").Value = -10349
$ws.Range("H20").Value = 10009
$ws.Range("J20").Value = 10009
$ws.Range("L20").Value = 10009
$ws.Range("N20").Value = -10549
$ws.Range("H37").Value = 6017
$ws.Range("I37").Value = 2034
$ws.Range("J37").Value = 10000
$ws.Range("K37").Value = 2034
$ws.Range("L37").Value = 10000
$ws.Range("M37").Value = -1761
$ws.Range("N37").Value = -10546
$ws.Range("H44").Value = 10500
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").ClearContents()
$ws.Range("H45").Value = 1768.75
$ws.Range("I45").Value = 1768.75
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 1768.75
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = -1391.75
$ws.Range("N45").ClearContents()
$ws.Range("H63").Value = 2228.5557
$ws.Range("I63").Value = 2228.5557
$ws.Range("K63").Value = 2228.5557
$ws.Range("M63").Value = -1542.5557
$ws.Range("H66").Value = 2228.5557
$ws.Range("I66").Value = 2228.5557
$ws.Range("K66").Value = 11142.7785
$ws.Range("M66").Value = -7710.7785
$ws.Range("H132").Value = 2607.7083
$ws.Range("I132").Value = 2189.1052
$ws.Range("J132").Value = 4198.4
$ws.Range("K132").Value = 6567.3156
$ws.Range("L132").Value = 12595.2
$ws.Range("M132").Value = -4037.3156
$ws.Range("N132").Value = -17655.2
$ws.Range("H133").Value = 27386.37
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 27386.37
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 27386.37
$ws.Range("M133").ClearContents()
$ws.Range("N133").Value = -32446.37

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 27814
$ws.Range("I82").Value = 25085.334
$ws.Range("K82").Value = 25085.334
$ws.Range("M82").Value = -24702.334
$ws.Range("N82").ClearContents()
$ws.Range("H85").Value = 27814
$ws.Range("I85").Value = 25085.334
$ws.Range("K85").Value = 25085.334
$ws.Range("M85").Value = -23759.334
$ws.Range("N85").ClearContents()
$ws.Range("H94").Value = 41667588
$ws.Range("I94").Value = 50000904
$ws.Range("J94").Value = 1010
$ws.Range("K94").Value = 50000904
$ws.Range("L94").Value = 1010
$ws.Range("M94").Value = -50000453
$ws.Range("N94").Value = -1912
$ws.Range("H134").Value = 5215.0386
$ws.Range("I134").Value = 1358
$ws.Range("J134").Value = 51499.5
$ws.Range("K134").Value = 4074
$ws.Range("L134").Value = 154498.5
$ws.Range("M134").Value = -1539
$ws.Range("N134").Value = -159568.5

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("M15").ClearContents()
$ws.Range("H22").Value = 600
$ws.Range("I22").Value = 600
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 600
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -250
$ws.Range("N22").ClearContents()
$ws.Range("H28").Value = 24000
$ws.Range("J28").Value = 24000
$ws.Range("L28").Value = 24000
$ws.Range("N28").Value = -24490
$ws.Range("H31").Value = 2278.3215
$ws.Range("I31").Value = 1190.7273
$ws.Range("K31").Value = 1190.7273
$ws.Range("M31").Value = -895.7273
$ws.Range("N31").ClearContents()
$ws.Range("H34").Value = 2278.3215
$ws.Range("I34").Value = 1190.7273
$ws.Range("K34").Value = 1190.7273
$ws.Range("M34").Value = -988.7273
$ws.Range("N34").ClearContents()
$ws.Range("H48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("N48").ClearContents()

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1387.1666
$ws.Range("I5").Value = 1803
$ws.Range("K5").Value = 5409
$ws.Range("M5").Value = -5297
$ws.Range("N5").ClearContents()
$ws.Range("H47").Value = 412.25
$ws.Range("I47").Value = 412.25
$ws.Range("K47").Value = 1236.75
$ws.Range("M47").Value = -805.75
$ws.Range("H59").Value = 4833
$ws.Range("J59").Value = 12000
$ws.Range("L59").Value = 36000
$ws.Range("N59").Value = -37080
$ws.Range("H135").Value = 1387.1666
$ws.Range("I135").Value = 1803
$ws.Range("K135").Value = 16227
$ws.Range("M135").Value = -13692
$ws.Range("N135").ClearContents()
$ws.Range("H137").Value = 2657.3333
$ws.Range("I137").Value = 1365
$ws.Range("J137").Value = 3691.2
$ws.Range("K137").Value = 4095
$ws.Range("L137").Value = 11073.6
$ws.Range("M137").Value = 1005
$ws.Range("N137").Value = -21273.6

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("H26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("N26").ClearContents()
$ws.Range("H50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("N50").ClearContents()

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 667.5
$ws.Range("I9").Value = 390
$ws.Range("K9").Value = 390
$ws.Range("M9").Value = -166
$ws.Range("N9").ClearContents()
$ws.Range("H22").Value = 1298.8823
$ws.Range("I22").Value = 1328
$ws.Range("J22").Value = 1257.2858
$ws.Range("K22").Value = 1328
$ws.Range("L22").Value = 1257.2858
$ws.Range("M22").Value = -1033
$ws.Range("N22").Value = -1847.2858
$ws.Range("H27").Value = 1298.8823
$ws.Range("I27").Value = 1328
$ws.Range("J27").Value = 1257.2858
$ws.Range("K27").Value = 1328
$ws.Range("L27").Value = 1257.2858
$ws.Range("M27").Value = -1221
$ws.Range("N27").Value = -1471.2858
$ws.Range("H40").Value = 3883.1667
$ws.Range("I40").Value = 3599.6667
$ws.Range("J40").Value = 4166.6665
$ws.Range("K40").Value = 3599.6667
$ws.Range("L40").Value = 4166.6665
$ws.Range("M40").Value = -3463.6667
$ws.Range("N40").Value = -4438.6665
$ws.Range("H56").Value = 7000
$ws.Range("I56").Value = 3666.6667
$ws.Range("K56").Value = 3666.6667
$ws.Range("M56").Value = -2975.6667
$ws.Range("N56").ClearContents()
$ws.Range("H122").Value = 17859868
$ws.Range("I122").Value = 31252568
$ws.Range("J122").Value = 2935.5
$ws.Range("K122").Value = 93757704
$ws.Range("L122").Value = 8806.5
$ws.Range("M122").Value = -93755254
$ws.Range("N122").Value = -13706.5
$ws.Range("H136").Value = 18734.666
$ws.Range("I136").Value = 34702.668
$ws.Range("J136").Value = 2766.6667
$ws.Range("K136").Value = 104108.004
$ws.Range("L136").Value = 8300.000100000001
$ws.Range("M136").Value = -101558.004
$ws.Range("N136").Value = -13400.0001

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H58").Value = 17999
$ws.Range("I58").Value = 17999
$ws.Range("K58").Value = 17999
$ws.Range("M58").Value = -17691
